$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: A2=0, B2=738
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 738

# Update row 3: A3=1, B3=264
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 264

# Remove rows 4 and 5 (old rows with A4=1/B4=204, A5=2/B5=173)
# so the sheet's used range shrinks back to A1:B3
$ws.Range("A4:B5").Delete()
